$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 88: fill in previously-empty cells with the "public testing" milestone ---
# (Note entry order matters for shared-string indices: the "Went out..." note
# was authored before the "Public testing Demo v0.3" event label.)
$ws.Range("D88").Value = "Went out in public for first time and got public to play test"
$ws.Range("B88").Value = "Public testing Demo v0.3"
$ws.Range("C88").Value = 44646

# A88 continues the same shared percent-complete formula used by A80:A104 (si="13")
$ws.Range("A88").Formula = "=MAX(0,(C88-`$C`$79)/(`$C`$103-`$C`$79))"

# --- New budget/finance figures in columns G/H ---
$ws.Range("G89").Value = 5487.1
$ws.Range("G103").Value = 390
$ws.Range("G104").Value = 1475
$ws.Range("G105").Formula = "=G104+G103"
$ws.Range("G93").Formula = "=G89-G105"
$ws.Range("G94").Value = 1733.33
$ws.Range("G95").Formula = "=G93/G94"
$ws.Range("H95").Formula = "=G95*B32"
$ws.Range("H96").Formula = "=TODAY()+H95"

# Match number formats seen in the target: G89/G93/G94/G105 use the "$"#,##0.00
# currency style (same as the burn-rate column, e.g. C8); H96 uses the plain
# date style used elsewhere in the sheet (e.g. B108). Copy/PasteSpecial (formats
# only) reuses the existing style slot instead of minting a fresh one.
$ws.Range("C8").Copy() | Out-Null
$ws.Range("G89").PasteSpecial(-4122) | Out-Null
$ws.Range("G93").PasteSpecial(-4122) | Out-Null
$ws.Range("G94").PasteSpecial(-4122) | Out-Null
$ws.Range("G105").PasteSpecial(-4122) | Out-Null

$ws.Range("B108").Copy() | Out-Null
$ws.Range("H96").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Sheet view: scroll position / selected cell moved ---
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("E39").Select()

$wb.Save()
